$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

# Row 2
Set-TextValue "D2" "90.761.15"
Set-TextValue "E2" "  -0.21%  "

# Row 3
Set-TextValue "D3" "3.114.53"
Set-TextValue "E3" "  -2.02%  "

# Row 4
Set-TextValue "D4" "0.995"
Set-TextValue "E4" "  -0.51%  "

# Row 5
Set-TextValue "D5" "231.69"
Set-TextValue "E5" "  +5.18%  "

# Row 6
Set-TextValue "D6" "626.43"
Set-TextValue "E6" "  +0.16%  "

# Row 7
Set-TextValue "D7" "1.11"
Set-TextValue "E7" "  +3.00%  "

# Row 8
Set-TextValue "D8" "0.365"
Set-TextValue "E8" "  -2.57%  "

# Row 9
Set-TextValue "D9" "0.999"
Set-TextValue "E9" "  -0.05%  "

# Row 10
Set-TextValue "D10" "3.112.17"
Set-TextValue "E10" "  -2.07%  "

# Row 11
Set-TextValue "D11" "0.724"
Set-TextValue "E11" "  -4.70%  "

# Row 12
Set-TextValue "E12" "  -0.85%  "

# Row 13
Set-TextValue "D13" "36.55"
Set-TextValue "E13" "  +3.34%  "

# Row 14
Set-TextValue "D14" "0.0000247"
Set-TextValue "E14" "  -1.72%  "

# Row 15
Set-TextValue "E15" "  -1.68%  "

# Row 16
Set-TextValue "D16" "90.645.30"
Set-TextValue "E16" "  +0.23%  "

# Row 17
Set-TextValue "D17" "3.693.36"
Set-TextValue "E17" "  -2.03%  "

# Row 18
Set-TextValue "D18" "3.114.31"
Set-TextValue "E18" "  -1.89%  "

# Row 19
Set-TextValue "D19" "3.81"
Set-TextValue "E19" "  +0.59%  "

# Row 20
Set-TextValue "D20" "14.14"
Set-TextValue "E20" "  -1.78%  "

# Row 21
Set-TextValue "D21" "0.0000209"
Set-TextValue "E21" "  -5.26%  "

# Row 22
Set-TextValue "D22" "441.33"
Set-TextValue "E22" "  -0.41%  "

# Row 23
Set-TextValue "D23" "5.56"
Set-TextValue "E23" "  +6.56%  "

# Row 24
Set-TextValue "D24" "8.90"
Set-TextValue "E24" "  -1.03%  "

# Row 25
Set-TextValue "D25" "5.88"
Set-TextValue "E25" "  -3.22%  "

# Row 26
Set-TextValue "D26" "89.13"
Set-TextValue "E26" "  +2.14%  "

# Row 27
Set-TextValue "D27" "12.34"
Set-TextValue "E27" "  -0.49%  "

# Row 28
Set-TextValue "D28" "3.315.53"
Set-TextValue "E28" "  -1.06%  "

# Row 29
Set-TextValue "D29" "0.999"
Set-TextValue "E29" "  -0.12%  "

# Row 30
Set-TextValue "D30" "9.46"
Set-TextValue "E30" "  +1.54%  "

# Row 31
Set-TextValue "D31" "0.159"
Set-TextValue "E31" "  -3.03%  "

# Row 32
Set-TextValue "D32" "0.200"
Set-TextValue "E32" "  +18.58%  "

# Row 33
Set-TextValue "D33" "26.47"
Set-TextValue "E33" "  +5.22%  "

# Row 34
Set-TextValue "D34" "0.894"
Set-TextValue "E34" "  -10.46%  "

# Row 35
Set-TextValue "E35" "  +3.59%  "

# Row 36
Set-TextValue "B36" "dogwifhat"
Set-TextValue "C36" "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue "D36" "3.78"
Set-TextValue "E36" "  +0.42%  "

# Row 37
Set-TextValue "B37" "Bittensor"
Set-TextValue "C37" "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue "D37" "510.82"
Set-TextValue "E37" "  -3.36%  "

# Row 38
Set-TextValue "B38" "RenderToken"
Set-TextValue "C38" "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
Set-TextValue "D38" "7.06"
Set-TextValue "E38" "  -0.31%  "

# Row 39
Set-TextValue "B39" "PancakeSwap"
Set-TextValue "C39" "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue "D39" "1.92"
Set-TextValue "E39" "  +0.42%  "

# Row 40
Set-TextValue "E40" "  -2.76%  "

# Row 41
Set-TextValue "B41" "Hedera"
Set-TextValue "C41" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D41" "0.0889"
Set-TextValue "E41" "  +5.70%  "

# Row 42
Set-TextValue "B42" "PolygonEcosystemToken"
Set-TextValue "C42" "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
Set-TextValue "D42" "0.411"
Set-TextValue "E42" "  -0.88%  "

# Row 43
Set-TextValue "D43" "22.19"
Set-TextValue "E43" "  -0.14%  "

# Row 45
Set-TextValue "D45" "3.39"
Set-TextValue "E45" "  +51.86%  "

# Row 46
Set-TextValue "D46" "1.91"
Set-TextValue "E46" "  -2.69%  "

# Row 47
Set-TextValue "D47" "151.06"
Set-TextValue "E47" "  +1.32%  "

# Row 48
Set-TextValue "B48" "ARBITRUM"
Set-TextValue "C48" "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue "D48" "0.689"
Set-TextValue "E48" "  +5.39%  "

# Row 49
Set-TextValue "B49" "OKB"
Set-TextValue "C49" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D49" "45.19"
Set-TextValue "E49" "  +2.06%  "

# Row 50
Set-TextValue "E50" "  -1.46%  "

# Row 51
Set-TextValue "D51" "4.46"
Set-TextValue "E51" "  +1.68%  "
